$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 40
$ws.Cells.Item(40, 1).Value = 43697
$ws.Cells.Item(40, 2).Value = 0.48333333333333334
$ws.Cells.Item(40, 3).Value = 72
$ws.Cells.Item(40, 4).Value = 21484

# New row 41
$ws.Cells.Item(41, 1).Value = 43697
$ws.Cells.Item(41, 2).Value = 0.50138888888888888
$ws.Cells.Item(41, 3).Value = 72
$ws.Cells.Item(41, 4).Value = 21449

# Copy number formatting from existing rows (A39 = date style, B5 = time style)
# so the new cells reuse the same cell styles instead of creating new ones.
$ws.Range("A39").Copy()
$ws.Range("A40:A41").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("B40:B41").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the view: scrolled position and active cell/selection to match the
# state Excel would save after scrolling down to show the newly added rows.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("A41").Select()
